$d = $word.ActiveDocument

# 1) The document used to open with a leading, otherwise-empty paragraph
#    that only held a manual page break (<w:br w:type="page"/>) so the
#    "Introduction" heading would always start on a new page. That lead-in
#    paragraph was removed, so delete it (and its paragraph mark) here.
$d.Paragraphs.Item(1).Range.Delete()

# 2) Because the page break used to force "Introduction" onto a new page,
#    its run still carries a <w:lastRenderedPageBreak/> rendering marker.
#    With the manual break gone that marker is stale, so rewrite the
#    paragraph's text (stripping the trailing paragraph-mark character
#    first) to force the run to be re-emitted without it.
$introPara = $d.Paragraphs.Item(1)
$introText = $introPara.Range.Text.TrimEnd([char]13)
$introPara.Range.Text = $introText
